# Updates cryptos list figures (price/volume) per GitHub Actions scrape refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text updates (coin name/link swap for rows 28-29, and volume/% + multi-dot
# prices that Excel never mistakes for numbers) - set directly via .Value.
$directUpdates = @(
    @("D2", "68.500.64"),
    @("E2", "  +0.73%  "),
    @("D3", "3.763.41"),
    @("E3", "  -0.49%  "),
    @("E4", "  +0.04%  "),
    @("E5", "  -0.48%  "),
    @("E6", "  -1.48%  "),
    @("D7", "3.762.04"),
    @("E7", "  -0.53%  "),
    @("E8", "  -0.06%  "),
    @("E9", "  -1.05%  "),
    @("E10", "  -2.59%  "),
    @("E11", "  -1.71%  "),
    @("E12", "  -1.08%  "),
    @("E13", "  -7.13%  "),
    @("E14", "  -1.57%  "),
    @("D15", "4.393.06"),
    @("E15", "  -0.50%  "),
    @("D16", "3.753.85"),
    @("E16", "  -0.83%  "),
    @("D17", "68.516.43"),
    @("E17", "  +0.89%  "),
    @("E18", "  -5.09%  "),
    @("E20", "  -3.26%  "),
    @("E21", "  +1.36%  "),
    @("E22", "  -0.41%  "),
    @("E23", "  -3.41%  "),
    @("E24", "  +0.40%  "),
    @("E25", "  -2.33%  "),
    @("E26", "  -2.91%  "),
    @("E27", "  -1.55%  "),
    @("B28", "RenderToken"),
    @("C28", "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"),
    @("E28", "  -4.67%  "),
    @("B29", "Dai"),
    @("C29", "https://coinranking.com/coin/MoTuySvg7+dai-dai"),
    @("E29", "  -0.09%  "),
    @("D30", "3.909.30"),
    @("E30", "  -0.52%  "),
    @("E31", "  -4.95%  "),
    @("E32", "  -3.70%  "),
    @("E33", "  -1.62%  "),
    @("E34", "  -3.19%  "),
    @("E35", "  -0.45%  "),
    @("D37", "3.715.33"),
    @("E37", "  -0.65%  "),
    @("E38", "  -3.79%  "),
    @("E39", "  -9.12%  "),
    @("E40", "  -0.30%  "),
    @("E41", "  +0.02%  "),
    @("E42", "  -1.04%  "),
    @("E43", "  +0.08%  "),
    @("E45", "  +9.45%  "),
    @("E46", "  -3.52%  "),
    @("E47", "  +2.62%  "),
    @("E48", "  -1.96%  "),
    @("E49", "  -2.43%  "),
    @("E50", "  +2.62%  "),
    @("E51", "  -3.56%  ")
)
foreach ($pair in $directUpdates) {
    $ws.Range($pair[0]).Value = $pair[1]
}

# Price updates that LOOK like plain numbers (e.g. "594.00", "0.0000260") must be
# forced to stay text - matching the sheet convention where the Price column is
# inline/shared string, not numeric. Prefix with an apostrophe (quote-prefix) so
# Excel stores the literal text, then clear the resulting quote-prefix style back
# to Normal so no stray number-format/style gets attached to the cell.
$quotedUpdates = @(
    @("D5", "594.00"),
    @("D6", "167.28"),
    @("D11", "6.44"),
    @("D13", "0.0000260"),
    @("D14", "36.06"),
    @("D18", "17.92"),
    @("D20", "6.99"),
    @("D22", "465.92"),
    @("D23", "0.696"),
    @("D24", "84.07"),
    @("D25", "0.0000147"),
    @("D26", "2.18"),
    @("D27", "11.96"),
    @("D28", "10.04"),
    @("D29", "1.00"),
    @("D31", "2.78"),
    @("D32", "7.32"),
    @("D33", "30.07"),
    @("D35", "9.20"),
    @("D38", "0.101"),
    @("D39", "3.41"),
    @("D42", "5.79"),
    @("D45", "44.17"),
    @("D46", "0.303"),
    @("D47", "46.83"),
    @("D49", "8.49"),
    @("D50", "145.47"),
    @("D51", "388.11")
)
foreach ($pair in $quotedUpdates) {
    $ws.Range($pair[0]).Value = "'" + $pair[1]
}
$ws.Range("D2:D51").Style = "Normal"
